# SCO-T1-A05-A06.pptx, slide 10 ("Modos de endereçamento"), shape "Rectangle 23"
# paragraph 3 ("Índice + pré-incremento ...").
#
# Original run (colored, Courier New, lang en-GB):   ", [R1 + #4]"   (comma, NBSP, text)
# becomes two runs:                                   ", [R1 + #"  +  "4]!"
# and the following run's leading closing-quote/comma ("", ") is split out into
# its own run, so the final run starts cleanly with "O registo R1 ...".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item(2)              # "Rectangle 23"
$tf = $shp.TextFrame
$tr = $tf.TextRange

$para3 = $tr.Paragraphs(3, 1)         # "Índice + pré-incremento: ..."
$start = $para3.Start

$nbsp = [char]0xA0
$rdquo = [char]0x201D                 # ” right double quotation mark

# --- Locate the run ", [R1 + #4]" inside paragraph 3 -----------------------
# (comma + NBSP + "[R1 + #4]") -- local (0-based) offset 164, length 11.
$runStart = $start + 164

# 1) Grow that run's text, adding the "!" -> ", [R1 + #4]!"
$run = $tr.Characters($runStart, 11)
$run.Text = "," + $nbsp + "[R1 + #4]!"

# 2) Split the run after "#" (9 chars in) so we get ", [R1 + #" / "4]!" as two
#    separate runs. Re-asserting a font attribute at its current value forces
#    the COM host to materialise a new run boundary without altering the
#    visible formatting.
$firstPart = $tr.Characters($runStart, 9)
$firstPart.Font.Italic = $false

# 3) Split the closing quote + comma + space ("”, ") off the front of the next
#    run, so it becomes its own run and the following run starts at "O registo".
$quotePart = $tr.Characters($runStart + 11 + 1, 3)
$quotePart.Text = $rdquo + ", "
